$wb = $excel.ActiveWorkbook

# Rename the first sheet (was inconsistently named after another table)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Data"

# Remove the stray blank formatted row in the "Column Information" sheet
$ws2 = $wb.Worksheets.Item("Column Information")
$ws2.Range("A3:C3").Clear()
